# This workbook is a weekly-updated price log ("Vega Monumental Concepción - Brócoli").
# A new week's worth of data (2 rows: calidad "Primera" and "Segunda") is inserted right
# before the existing data block (originally starting at row 535), pushing all the
# existing rows down by two. The two new rows are populated with the latest week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block (row 535), shifting
# everything currently at/after row 535 down by two rows.
$ws.Range("535:536").Insert()

# --- New row 535: Calidad "Primera" ---
$ws.Cells.Item(535, 1).Value = 11
$ws.Cells.Item(535, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(535, 3).Value = "Bíobío"
$ws.Cells.Item(535, 4).Value = 45212
$ws.Cells.Item(535, 5).Value = 8
$ws.Cells.Item(535, 6).Value = 100112023
$ws.Cells.Item(535, 7).Value = "Brócoli"
$ws.Cells.Item(535, 8).Value = "Sin especificar"
$ws.Cells.Item(535, 9).Value = "Primera"
$ws.Cells.Item(535, 10).Value = 2000
$ws.Cells.Item(535, 11).Value = 800
$ws.Cells.Item(535, 12).Value = 900
$ws.Cells.Item(535, 13).Value = 850
$ws.Cells.Item(535, 14).Value = "`$/unidad"
$ws.Cells.Item(535, 15).Value = "Región Metropolitana"
$ws.Cells.Item(535, 16).Value = 850
$ws.Cells.Item(535, 17).Value = 1
$ws.Cells.Item(535, 18).Value = "Hortaliza"

# --- New row 536: Calidad "Segunda" ---
$ws.Cells.Item(536, 1).Value = 11
$ws.Cells.Item(536, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(536, 3).Value = "Bíobío"
$ws.Cells.Item(536, 4).Value = 45212
$ws.Cells.Item(536, 5).Value = 8
$ws.Cells.Item(536, 6).Value = 100112023
$ws.Cells.Item(536, 7).Value = "Brócoli"
$ws.Cells.Item(536, 8).Value = "Sin especificar"
$ws.Cells.Item(536, 9).Value = "Segunda"
$ws.Cells.Item(536, 10).Value = 1000
$ws.Cells.Item(536, 11).Value = 700
$ws.Cells.Item(536, 12).Value = 700
$ws.Cells.Item(536, 13).Value = 700
$ws.Cells.Item(536, 14).Value = "`$/unidad"
$ws.Cells.Item(536, 15).Value = "Región Metropolitana"
$ws.Cells.Item(536, 16).Value = 700
$ws.Cells.Item(536, 17).Value = 1
$ws.Cells.Item(536, 18).Value = "Hortaliza"

# Make sure the date cells use the same date/time number format as the rest of
# column D (the Insert() above should already have carried this down, but set
# it explicitly to be safe).
$ws.Cells.Item(535, 4).NumberFormat = $ws.Cells.Item(537, 4).NumberFormat
$ws.Cells.Item(536, 4).NumberFormat = $ws.Cells.Item(537, 4).NumberFormat
